$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 392, shifting rows 392:451 down to 393:452.
$ws.Rows("392:392").Insert()

# Populate the new row 392 with the new data point.
$ws.Range("A392").Value = 9
$ws.Range("B392").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C392").Value = "Metropolitana"
$ws.Range("D392").Value = 45180
$ws.Range("E392").Value = 13
$ws.Range("F392").Value = 100112043
$ws.Range("G392").Value = "Pepino ensalada"
$ws.Range("H392").Value = "Sin especificar"
$ws.Range("I392").Value = "Primera"
$ws.Range("J392").Value = 70
$ws.Range("K392").Value = 11000
$ws.Range("L392").Value = 13000
$ws.Range("M392").Value = 12000
$ws.Range("N392").Value = "$/caja 60 unidades"
$ws.Range("O392").Value = "Región de Arica y Parinacota"
$ws.Range("P392").Value = 200
$ws.Range("Q392").Value = 60
$ws.Range("R392").Value = "Hortaliza"
